$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-02-12 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-02-13 Thursday", 2) | Out-Null
$d.Content.Find.Execute("35-1=", $true, $false, $false, $false, $false, $true, 1, $false, "8-3=", 2) | Out-Null
$d.Content.Find.Execute("29-10=", $true, $false, $false, $false, $false, $true, 1, $false, "57+42=", 2) | Out-Null
$d.Content.Find.Execute("77-5=", $true, $false, $false, $false, $false, $true, 1, $false, "16-11=", 2) | Out-Null
$d.Content.Find.Execute("18-9=", $true, $false, $false, $false, $false, $true, 1, $false, "10+13=", 2) | Out-Null
$d.Content.Find.Execute("46+17=", $true, $false, $false, $false, $false, $true, 1, $false, "88-17=", 2) | Out-Null
$d.Content.Find.Execute("87-50=", $true, $false, $false, $false, $false, $true, 1, $false, "9+16=", 2) | Out-Null
$d.Content.Find.Execute("27+26=", $true, $false, $false, $false, $false, $true, 1, $false, "12+10=", 2) | Out-Null
$d.Content.Find.Execute("0+56=", $true, $false, $false, $false, $false, $true, 1, $false, "79-61=", 2) | Out-Null
$d.Content.Find.Execute("61-41=", $true, $false, $false, $false, $false, $true, 1, $false, "71-32=", 2) | Out-Null
$d.Content.Find.Execute("15+62=", $true, $false, $false, $false, $false, $true, 1, $false, "36+16=", 2) | Out-Null
$d.Content.Find.Execute("12-0=", $true, $false, $false, $false, $false, $true, 1, $false, "63-21=", 2) | Out-Null
$d.Content.Find.Execute("46+20=", $true, $false, $false, $false, $false, $true, 1, $false, "57+34=", 2) | Out-Null
$d.Content.Find.Execute("36+51=", $true, $false, $false, $false, $false, $true, 1, $false, "63-48=", 2) | Out-Null
$d.Content.Find.Execute("96-74=", $true, $false, $false, $false, $false, $true, 1, $false, "35-35=", 2) | Out-Null
$d.Content.Find.Execute("12+2=", $true, $false, $false, $false, $false, $true, 1, $false, "11-2=", 2) | Out-Null
$d.Content.Find.Execute("61-30=", $true, $false, $false, $false, $false, $true, 1, $false, "99-52=", 2) | Out-Null
$d.Content.Find.Execute("44-32=", $true, $false, $false, $false, $false, $true, 1, $false, "1+13=", 2) | Out-Null
$d.Content.Find.Execute("24-14=", $true, $false, $false, $false, $false, $true, 1, $false, "13+33=", 2) | Out-Null
$d.Content.Find.Execute("26+65=", $true, $false, $false, $false, $false, $true, 1, $false, "40-11=", 2) | Out-Null
$d.Content.Find.Execute("40-24=", $true, $false, $false, $false, $false, $true, 1, $false, "45+6=", 2) | Out-Null
$d.Content.Find.Execute("88-72=", $true, $false, $false, $false, $false, $true, 1, $false, "50+19=", 2) | Out-Null
$d.Content.Find.Execute("78-62=", $true, $false, $false, $false, $false, $true, 1, $false, "38+33=", 2) | Out-Null
$d.Content.Find.Execute("40+59=", $true, $false, $false, $false, $false, $true, 1, $false, "52-46=", 2) | Out-Null
$d.Content.Find.Execute("61-2=", $true, $false, $false, $false, $false, $true, 1, $false, "23+55=", 2) | Out-Null
$d.Content.Find.Execute("46+5=", $true, $false, $false, $false, $false, $true, 1, $false, "48+35=", 2) | Out-Null
$d.Content.Find.Execute("36+40=", $true, $false, $false, $false, $false, $true, 1, $false, "54+22=", 2) | Out-Null
$d.Content.Find.Execute("83+3=", $true, $false, $false, $false, $false, $true, 1, $false, "34-4=", 2) | Out-Null
$d.Content.Find.Execute("12+33=", $true, $false, $false, $false, $false, $true, 1, $false, "24+44=", 2) | Out-Null
$d.Content.Find.Execute("57-41=", $true, $false, $false, $false, $false, $true, 1, $false, "52-43=", 2) | Out-Null
$d.Content.Find.Execute("57+37=", $true, $false, $false, $false, $false, $true, 1, $false, "95-68=", 2) | Out-Null
$d.Content.Find.Execute("23+72=", $true, $false, $false, $false, $false, $true, 1, $false, "54+13=", 2) | Out-Null
$d.Content.Find.Execute("29-6=", $true, $false, $false, $false, $false, $true, 1, $false, "70-62=", 2) | Out-Null
$d.Content.Find.Execute("95-15=", $true, $false, $false, $false, $false, $true, 1, $false, "5+74=", 2) | Out-Null
$d.Content.Find.Execute("58+19=", $true, $false, $false, $false, $false, $true, 1, $false, "86+10=", 2) | Out-Null
$d.Content.Find.Execute("33+52=", $true, $false, $false, $false, $false, $true, 1, $false, "55-12=", 2) | Out-Null
$d.Content.Find.Execute("92-69=", $true, $false, $false, $false, $false, $true, 1, $false, "47-45=", 2) | Out-Null
$d.Content.Find.Execute("59-43=", $true, $false, $false, $false, $false, $true, 1, $false, "91-13=", 2) | Out-Null
$d.Content.Find.Execute("38+60=", $true, $false, $false, $false, $false, $true, 1, $false, "41-4=", 2) | Out-Null
$d.Content.Find.Execute("95-77=", $true, $false, $false, $false, $false, $true, 1, $false, "71-54=", 2) | Out-Null
$d.Content.Find.Execute("37+48=", $true, $false, $false, $false, $false, $true, 1, $false, "61-50=", 2) | Out-Null
$d.Content.Find.Execute("22+52=", $true, $false, $false, $false, $false, $true, 1, $false, "9+73=", 2) | Out-Null
$d.Content.Find.Execute("15+76=", $true, $false, $false, $false, $false, $true, 1, $false, "88-59=", 2) | Out-Null
$d.Content.Find.Execute("71-29=", $true, $false, $false, $false, $false, $true, 1, $false, "60-32=", 2) | Out-Null
$d.Content.Find.Execute("20+75=", $true, $false, $false, $false, $false, $true, 1, $false, "48+1=", 2) | Out-Null
$d.Content.Find.Execute("4+19=", $true, $false, $false, $false, $false, $true, 1, $false, "4+74=", 2) | Out-Null
$d.Content.Find.Execute("45+34=", $true, $false, $false, $false, $false, $true, 1, $false, "77+22=", 2) | Out-Null
$d.Content.Find.Execute("30+64=", $true, $false, $false, $false, $false, $true, 1, $false, "91-77=", 2) | Out-Null
$d.Content.Find.Execute("27+14=", $true, $false, $false, $false, $false, $true, 1, $false, "31-9=", 2) | Out-Null
$d.Content.Find.Execute("63+22=", $true, $false, $false, $false, $false, $true, 1, $false, "26+12=", 2) | Out-Null
$d.Content.Find.Execute("7+74=", $true, $false, $false, $false, $false, $true, 1, $false, "97-93=", 2) | Out-Null
$d.Content.Find.Execute("63-9=", $true, $false, $false, $false, $false, $true, 1, $false, "67+26=", 2) | Out-Null
$d.Content.Find.Execute("49+26=", $true, $false, $false, $false, $false, $true, 1, $false, "74-61=", 2) | Out-Null
$d.Content.Find.Execute("92-39=", $true, $false, $false, $false, $false, $true, 1, $false, "70-64=", 2) | Out-Null
$d.Content.Find.Execute("3-2=", $true, $false, $false, $false, $false, $true, 1, $false, "92+7=", 2) | Out-Null
$d.Content.Find.Execute("6+11=", $true, $false, $false, $false, $false, $true, 1, $false, "58-11=", 2) | Out-Null
$d.Content.Find.Execute("42+54=", $true, $false, $false, $false, $false, $true, 1, $false, "61-47=", 2) | Out-Null
$d.Content.Find.Execute("41+40=", $true, $false, $false, $false, $false, $true, 1, $false, "76-55=", 2) | Out-Null
$d.Content.Find.Execute("55-48=", $true, $false, $false, $false, $false, $true, 1, $false, "27+7=", 2) | Out-Null
$d.Content.Find.Execute("42+57=", $true, $false, $false, $false, $false, $true, 1, $false, "67+18=", 2) | Out-Null
$d.Content.Find.Execute("16+2=", $true, $false, $false, $false, $false, $true, 1, $false, "76-72=", 2) | Out-Null
$d.Content.Find.Execute("18+39=", $true, $false, $false, $false, $false, $true, 1, $false, "34+42=", 2) | Out-Null
$d.Content.Find.Execute("80-49=", $true, $false, $false, $false, $false, $true, 1, $false, "85-8=", 2) | Out-Null
$d.Content.Find.Execute("80-75=", $true, $false, $false, $false, $false, $true, 1, $false, "41-28=", 2) | Out-Null
$d.Content.Find.Execute("90-64=", $true, $false, $false, $false, $false, $true, 1, $false, "39+59=", 2) | Out-Null
$d.Content.Find.Execute("47-42=", $true, $false, $false, $false, $false, $true, 1, $false, "98-26=", 2) | Out-Null
$d.Content.Find.Execute("84-38=", $true, $false, $false, $false, $false, $true, 1, $false, "39+36=", 2) | Out-Null
$d.Content.Find.Execute("80-48=", $true, $false, $false, $false, $false, $true, 1, $false, "95-58=", 2) | Out-Null
$d.Content.Find.Execute("62+5=", $true, $false, $false, $false, $false, $true, 1, $false, "69-4=", 2) | Out-Null
$d.Content.Find.Execute("51-48=", $true, $false, $false, $false, $false, $true, 1, $false, "6+0=", 2) | Out-Null
$d.Content.Find.Execute("84-59=", $true, $false, $false, $false, $false, $true, 1, $false, "24+27=", 2) | Out-Null
$d.Content.Find.Execute("6+60=", $true, $false, $false, $false, $false, $true, 1, $false, "64-40=", 2) | Out-Null
$d.Content.Find.Execute("49-48=", $true, $false, $false, $false, $false, $true, 1, $false, "57-43=", 2) | Out-Null
$d.Content.Find.Execute("8+78=", $true, $false, $false, $false, $false, $true, 1, $false, "29+12=", 2) | Out-Null
$d.Content.Find.Execute("16+32=", $true, $false, $false, $false, $false, $true, 1, $false, "81+15=", 2) | Out-Null
$d.Content.Find.Execute("82-13=", $true, $false, $false, $false, $false, $true, 1, $false, "54-30=", 2) | Out-Null
$d.Content.Find.Execute("96-61=", $true, $false, $false, $false, $false, $true, 1, $false, "19+75=", 2) | Out-Null
$d.Content.Find.Execute("26+44=", $true, $false, $false, $false, $false, $true, 1, $false, "1+46=", 2) | Out-Null
$d.Content.Find.Execute("79-46=", $true, $false, $false, $false, $false, $true, 1, $false, "5+16=", 2) | Out-Null
$d.Content.Find.Execute("71-28=", $true, $false, $false, $false, $false, $true, 1, $false, "84-3=", 2) | Out-Null
$d.Content.Find.Execute("28-6=", $true, $false, $false, $false, $false, $true, 1, $false, "77+13=", 2) | Out-Null
$d.Content.Find.Execute("8+66=", $true, $false, $false, $false, $false, $true, 1, $false, "88-34=", 2) | Out-Null
$d.Content.Find.Execute("98-31=", $true, $false, $false, $false, $false, $true, 1, $false, "30-2=", 2) | Out-Null
$d.Content.Find.Execute("90-25=", $true, $false, $false, $false, $false, $true, 1, $false, "17-12=", 2) | Out-Null
$d.Content.Find.Execute("62-30=", $true, $false, $false, $false, $false, $true, 1, $false, "96-35=", 2) | Out-Null
$d.Content.Find.Execute("47-43=", $true, $false, $false, $false, $false, $true, 1, $false, "13+62=", 2) | Out-Null
$d.Content.Find.Execute("12+76=", $true, $false, $false, $false, $false, $true, 1, $false, "21+65=", 2) | Out-Null
$d.Content.Find.Execute("9+64=", $true, $false, $false, $false, $false, $true, 1, $false, "80+5=", 2) | Out-Null
$d.Content.Find.Execute("90-88=", $true, $false, $false, $false, $false, $true, 1, $false, "84+9=", 2) | Out-Null
$d.Content.Find.Execute("44-30=", $true, $false, $false, $false, $false, $true, 1, $false, "30-24=", 2) | Out-Null
$d.Content.Find.Execute("87+12=", $true, $false, $false, $false, $false, $true, 1, $false, "19+1=", 2) | Out-Null
$d.Content.Find.Execute("68-14=", $true, $false, $false, $false, $false, $true, 1, $false, "92-85=", 2) | Out-Null
$d.Content.Find.Execute("33-26=", $true, $false, $false, $false, $false, $true, 1, $false, "99-71=", 2) | Out-Null
$d.Content.Find.Execute("83-35=", $true, $false, $false, $false, $false, $true, 1, $false, "53-51=", 2) | Out-Null
$d.Content.Find.Execute("19-4=", $true, $false, $false, $false, $false, $true, 1, $false, "76-47=", 2) | Out-Null
$d.Content.Find.Execute("75-11=", $true, $false, $false, $false, $false, $true, 1, $false, "58-49=", 2) | Out-Null
$d.Content.Find.Execute("42+13=", $true, $false, $false, $false, $false, $true, 1, $false, "50+25=", 2) | Out-Null
$d.Content.Find.Execute("94-80=", $true, $false, $false, $false, $false, $true, 1, $false, "22+76=", 2) | Out-Null
$d.Content.Find.Execute("82-43=", $true, $false, $false, $false, $false, $true, 1, $false, "19+57=", 2) | Out-Null
$d.Content.Find.Execute("69-30=", $true, $false, $false, $false, $false, $true, 1, $false, "49-25=", 2) | Out-Null
$d.Content.Find.Execute("93+1=", $true, $false, $false, $false, $false, $true, 1, $false, "59-40=", 2) | Out-Null
